# Updates cryptos list values (Price / Volume(1h) columns) to the latest
# snapshot, and swaps the Cardano/TRON rows (12 <-> 13) back to their
# refreshed ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.948.24"
$ws.Range("E2").Value = "  +5.25%  "

$ws.Range("D3").Value = "2.723.57"
$ws.Range("E3").Value = "  +2.80%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D5").Value = "576.20"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D6").Value = "154.34"
$ws.Range("E6").Value = "  +6.07%  "

$ws.Range("D7").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "2.747.25"
$ws.Range("E9").Value = "  +2.93%  "

$ws.Range("D10").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("D11").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +5.11%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D12").Value = "0.162"
$ws.Range("E12").Value = "  +5.41%  "

$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D13").Value = "0.389"
$ws.Range("E13").Value = "  +3.54%  "

$ws.Range("D14").Value = "3.209.42"
$ws.Range("E14").Value = "  +2.90%  "

$ws.Range("D15").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D15").Value = "26.30"
$ws.Range("E15").Value = "  +2.36%  "

$ws.Range("D16").Value = "63.814.46"
$ws.Range("E16").Value = "  +5.06%  "

$ws.Range("E17").Value = "  +6.18%  "

$ws.Range("D18").Value = "2.742.26"
$ws.Range("E18").Value = "  +3.11%  "

$ws.Range("D19").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").Value = "  +3.17%  "

$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D21").Value = "360.34"
$ws.Range("E21").Value = "  +2.75%  "

$ws.Range("D22").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D22").Value = "6.93"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D24").Value = "0.533"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D25").Value = "66.03"
$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("E26").Value = "  +4.41%  "

$ws.Range("D27").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D27").Value = "8.50"
$ws.Range("E27").Value = "  +4.53%  "

$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").Value = "0.0₃0906"
$ws.Range("E29").Value = "  +11.83%  "

$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("D31").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  +6.30%  "

$ws.Range("D32").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D32").Value = "171.48"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("E33").Value = "  +13.21%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D35").Value = "20.47"
$ws.Range("E35").Value = "  +3.79%  "

$ws.Range("D36").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D36").Value = "4.77"
$ws.Range("E36").Value = "  +6.87%  "

$ws.Range("E37").Value = "  +9.14%  "

$ws.Range("E38").Value = "  +9.76%  "

$ws.Range("E39").Value = "  +13.17%  "

$ws.Range("D40").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D40").Value = "344.56"
$ws.Range("E40").Value = "  +2.80%  "

$ws.Range("E41").Value = "  +5.42%  "

$ws.Range("D42").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D42").Value = "39.26"
$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("D43").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  +7.36%  "

$ws.Range("D44").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D44").Value = "21.67"
$ws.Range("E44").Value = "  +6.24%  "

$ws.Range("D45").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D45").Value = "21.75"
$ws.Range("E45").Value = "  +3.90%  "

$ws.Range("E46").Value = "  +5.02%  "

$ws.Range("D47").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D47").Value = "139.02"
$ws.Range("E47").Value = "  +3.89%  "

$ws.Range("D48").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D48").Value = "0.643"
$ws.Range("E48").Value = "  +4.37%  "

$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"  # keep as text, not a number
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.01%  "
